$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 197; this shifts the existing rows 197-248 down to 198-249.
$ws.Rows.Item(197).Insert()

# Populate the newly inserted row 197 with the new record.
$ws.Cells.Item(197, 1).Value = 5
$ws.Cells.Item(197, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(197, 3).Value = "Maule"
$ws.Cells.Item(197, 4).Value = 44754
$ws.Cells.Item(197, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(197, 5).Value = 7
$ws.Cells.Item(197, 6).Value = 100112024
$ws.Cells.Item(197, 7).Value = "Choclo"
$ws.Cells.Item(197, 8).Value = "Dulce o Americano"
$ws.Cells.Item(197, 9).Value = "Primera"
$ws.Cells.Item(197, 10).Value = 100
$ws.Cells.Item(197, 11).Value = 38000
$ws.Cells.Item(197, 12).Value = 38000
$ws.Cells.Item(197, 13).Value = 38000
$ws.Cells.Item(197, 14).Value = "$/malla 70 unidades"
$ws.Cells.Item(197, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(197, 16).Value = 543
$ws.Cells.Item(197, 17).Value = 70
$ws.Cells.Item(197, 18).Value = "Hortaliza"
